$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date (row 2 / b9c71cc8...) refreshed
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 23:03:59"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2) refreshed
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-25 23:03:55"
$wsZhCn.Range("K2").Value = "2016-08-25 23:04:16"

# "de-de" sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2) refreshed
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-25 23:03:59"
$wsDeDe.Range("K2").Value = "2016-08-25 23:04:23"
